$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe so values that look numeric (e.g. "4.30") are
# stored verbatim as text instead of being auto-coerced to a Number by Excels
# smart-entry parser, then reset the cell style so no stray number-format /
# quote-prefix style gets attached (matches the source cells, which carry no
# explicit style).

$ws.Range('D2').Value = "'34.603.27"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.28%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.797.31"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.90%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.06%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'227.19"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.44%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.559"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +2.17%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.08%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'32.95"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +3.74%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.83%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.0697"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.99%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.43%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.91%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.814.50"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.01%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'11.08"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.75%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.639"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.53%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'34.580.81"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.36%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'4.30"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.98%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'68.97"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.47%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'247.67"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.25%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0803"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.30%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'11.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +2.92%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.12%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'4.18"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.93%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'168.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +3.36%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.42%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +1.77%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'16.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.84%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +2.25%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.97%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'4.11"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +11.88%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.12%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.82%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.95%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.78%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.429.34"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.94%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.58"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +7.02%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.98%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.25%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.63%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'85.73"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +6.85%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'2.42"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.01%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.27%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D44').Value = "'13.77"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.44%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0528"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +3.55%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'6.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.95%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.73%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.79%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'106.09"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +1.45%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.03%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.0₆0129"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -4.98%  "
$ws.Range('E51').Style = 'Normal'
